$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashify")

# Update the trade-in figures in column T (last period) for several product rows
$ws.Range("T2").Value = 72
$ws.Range("T4").Value = 425
$ws.Range("T5").Value = 23
$ws.Range("T6").Value = 0
$ws.Range("T7").Value = 7

# Update the active selection/cell shown in the sheet view
$ws.Range("W6").Select() | Out-Null
